$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 46 (toy-spam confidence bucket rebinned with min 5 threshold)
$ws.Rows.Item(46).Delete()

# Clear cells no longer populated
$ws.Range("J19,K19,L19,M19,N19,O19,P19,Q19").ClearContents()

# Update changed cell values
$ws.Range("B3").Value = 0.9347826086956522
$ws.Range("C3").Value = 43
$ws.Range("D3").Value = 43
$ws.Range("H3").Value = 3
$ws.Range("K3").Value = 0.8392857142857143
$ws.Range("L3").Value = 47
$ws.Range("M3").Value = 47
$ws.Range("Q3").Value = 9
$ws.Range("B4").Value = 0.8636363636363636
$ws.Range("C4").Value = 38
$ws.Range("D4").Value = 38
$ws.Range("H4").Value = 6
$ws.Range("K4").Value = 0.8153846153846154
$ws.Range("L4").Value = 53
$ws.Range("M4").Value = 53
$ws.Range("Q4").Value = 12
$ws.Range("B5").Value = 0.75
$ws.Range("C5").Value = 48
$ws.Range("D5").Value = 48
$ws.Range("H5").Value = 16
$ws.Range("J5").Value = 'favorite'
$ws.Range("K5").Value = 0.6989247311827957
$ws.Range("L5").Value = 65
$ws.Range("M5").Value = 65
$ws.Range("Q5").Value = 28
$ws.Range("B6").Value = 0.7311827956989247
$ws.Range("C6").Value = 136
$ws.Range("D6").Value = 136
$ws.Range("H6").Value = 50
$ws.Range("J6").Value = 'classic'
$ws.Range("K6").Value = 0.6226415094339622
$ws.Range("L6").Value = 33
$ws.Range("M6").Value = 33
$ws.Range("Q6").Value = 20
$ws.Range("A7").Value = 'broke'
$ws.Range("B7").Value = 0.7281553398058253
$ws.Range("C7").Value = 150
$ws.Range("D7").Value = 150
$ws.Range("H7").Value = 56
$ws.Range("J7").Value = 'excellent'
$ws.Range("K7").Value = 0.53125
$ws.Range("L7").Value = 34
$ws.Range("M7").Value = 34
$ws.Range("Q7").Value = 30
$ws.Range("A8").Value = 'returned'
$ws.Range("B8").Value = 0.7105263157894737
$ws.Range("C8").Value = 27
$ws.Range("D8").Value = 27
$ws.Range("H8").Value = 11
$ws.Range("J8").Value = 'thank'
$ws.Range("K8").Value = 0.391304347826087
$ws.Range("L8").Value = 27
$ws.Range("M8").Value = 27
$ws.Range("Q8").Value = 42
$ws.Range("A9").Value = 'poor'
$ws.Range("B9").Value = 0.6901408450704225
$ws.Range("C9").Value = 49
$ws.Range("D9").Value = 49
$ws.Range("H9").Value = 22
$ws.Range("J9").Value = 'great'
$ws.Range("K9").Value = 0.340983606557377
$ws.Range("L9").Value = 416
$ws.Range("M9").Value = 416
$ws.Range("Q9").Value = 804
$ws.Range("A10").Value = 'waste'
$ws.Range("B10").Value = 0.6216216216216216
$ws.Range("C10").Value = 92
$ws.Range("D10").Value = 92
$ws.Range("H10").Value = 56
$ws.Range("J10").Value = 'love'
$ws.Range("K10").Value = 0.2801724137931034
$ws.Range("L10").Value = 195
$ws.Range("M10").Value = 196
$ws.Range("N10").Value = 0.99
$ws.Range("O10").Value = 0.01000000000000001
$ws.Range("P10").Value = $true
$ws.Range("Q10").Value = 501
$ws.Range("A11").Value = 'water'
$ws.Range("B11").Value = 0.5714285714285714
$ws.Range("C11").Value = 24
$ws.Range("D11").Value = 24
$ws.Range("H11").Value = 18
$ws.Range("J11").Value = 'loves'
$ws.Range("K11").Value = 0.2572614107883817
$ws.Range("L11").Value = 124
$ws.Range("M11").Value = 124
$ws.Range("Q11").Value = 358
$ws.Range("A12").Value = 'smaller'
$ws.Range("B12").Value = 0.5546218487394958
$ws.Range("C12").Value = 66
$ws.Range("D12").Value = 66
$ws.Range("H12").Value = 53
$ws.Range("J12").Value = 'best'
$ws.Range("K12").Value = 0.2333333333333333
$ws.Range("L12").Value = 28
$ws.Range("M12").Value = 28
$ws.Range("Q12").Value = 92
$ws.Range("A13").Value = 'junk'
$ws.Range("B13").Value = 0.5272727272727272
$ws.Range("C13").Value = 29
$ws.Range("D13").Value = 29
$ws.Range("H13").Value = 26
$ws.Range("J13").Value = 'loved'
$ws.Range("K13").Value = 0.1896024464831804
$ws.Range("L13").Value = 62
$ws.Range("M13").Value = 62
$ws.Range("Q13").Value = 265
$ws.Range("A14").Value = 'broken'
$ws.Range("B14").Value = 0.4939759036144578
$ws.Range("C14").Value = 41
$ws.Range("D14").Value = 41
$ws.Range("H14").Value = 42
$ws.Range("K14").Value = 0.1807228915662651
$ws.Range("L14").Value = 30
$ws.Range("M14").Value = 30
$ws.Range("Q14").Value = 136
$ws.Range("A15").Value = 'guess'
$ws.Range("B15").Value = 0.4814814814814815
$ws.Range("C15").Value = 26
$ws.Range("D15").Value = 26
$ws.Range("H15").Value = 28
$ws.Range("J15").Value = 'friends'
$ws.Range("K15").Value = 0.164021164021164
$ws.Range("L15").Value = 31
$ws.Range("M15").Value = 31
$ws.Range("Q15").Value = 158
$ws.Range("A16").Value = 'small'
$ws.Range("B16").Value = 0.472463768115942
$ws.Range("C16").Value = 163
$ws.Range("D16").Value = 163
$ws.Range("H16").Value = 182
$ws.Range("J16").Value = 'christmas'
$ws.Range("K16").Value = 0.09236947791164658
$ws.Range("L16").Value = 23
$ws.Range("M16").Value = 23
$ws.Range("Q16").Value = 226
$ws.Range("A17").Value = 'instead'
$ws.Range("B17").Value = 0.4583333333333333
$ws.Range("C17").Value = 22
$ws.Range("D17").Value = 22
$ws.Range("H17").Value = 26
$ws.Range("J17").Value = 'fun'
$ws.Range("K17").Value = 0.08764241893076249
$ws.Range("L17").Value = 100
$ws.Range("M17").Value = 100
$ws.Range("Q17").Value = 1041
$ws.Range("A18").Value = 'paint'
$ws.Range("B18").Value = 0.4444444444444444
$ws.Range("C18").Value = 28
$ws.Range("D18").Value = 28
$ws.Range("H18").Value = 35
$ws.Range("J18").Value = 'game'
$ws.Range("K18").Value = 0.03634003893575601
$ws.Range("L18").Value = 56
$ws.Range("M18").Value = 56
$ws.Range("Q18").Value = 1485
$ws.Range("A19").Value = 'apart'
$ws.Range("B19").Value = 0.4210526315789473
$ws.Range("C19").Value = 40
$ws.Range("D19").Value = 40
$ws.Range("H19").Value = 55
$ws.Range("A20").Value = 'plastic'
$ws.Range("B20").Value = 0.4094488188976378
$ws.Range("C20").Value = 52
$ws.Range("D20").Value = 52
$ws.Range("H20").Value = 75
$ws.Range("A21").Value = 'di'
$ws.Range("B21").Value = 0.34375
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 22
$ws.Range("H21").Value = 42
$ws.Range("A22").Value = 'difficult'
$ws.Range("B22").Value = 0.3146067415730337
$ws.Range("C22").Value = 28
$ws.Range("D22").Value = 28
$ws.Range("H22").Value = 61
$ws.Range("A23").Value = 'thought'
$ws.Range("B23").Value = 0.301980198019802
$ws.Range("C23").Value = 61
$ws.Range("D23").Value = 61
$ws.Range("H23").Value = 141
$ws.Range("A24").Value = 'ok'
$ws.Range("B24").Value = 0.28125
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 36
$ws.Range("H24").Value = 92
$ws.Range("A25").Value = 'cheap'
$ws.Range("B25").Value = 0.2796208530805687
$ws.Range("C25").Value = 59
$ws.Range("D25").Value = 59
$ws.Range("H25").Value = 152
$ws.Range("A26").Value = 'size'
$ws.Range("B26").Value = 0.2371134020618557
$ws.Range("C26").Value = 46
$ws.Range("D26").Value = 46
$ws.Range("H26").Value = 148
$ws.Range("A27").Value = 'bit'
$ws.Range("B27").Value = 0.2346938775510204
$ws.Range("C27").Value = 23
$ws.Range("D27").Value = 23
$ws.Range("H27").Value = 75
$ws.Range("A28").Value = 'though'
$ws.Range("B28").Value = 0.2307692307692308
$ws.Range("C28").Value = 27
$ws.Range("D28").Value = 27
$ws.Range("H28").Value = 90
$ws.Range("A29").Value = 'item'
$ws.Range("B29").Value = 0.2028985507246377
$ws.Range("C29").Value = 56
$ws.Range("D29").Value = 56
$ws.Range("H29").Value = 220
$ws.Range("A30").Value = 'hard'
$ws.Range("B30").Value = 0.195
$ws.Range("C30").Value = 39
$ws.Range("D30").Value = 39
$ws.Range("H30").Value = 161
$ws.Range("A31").Value = 'money'
$ws.Range("B31").Value = 0.1930379746835443
$ws.Range("C31").Value = 61
$ws.Range("D31").Value = 61
$ws.Range("H31").Value = 255
$ws.Range("A32").Value = '1'
$ws.Range("B32").Value = 0.1864406779661017
$ws.Range("C32").Value = 22
$ws.Range("D32").Value = 22
$ws.Range("H32").Value = 96
$ws.Range("A33").Value = 'would'
$ws.Range("B33").Value = 0.1824925816023739
$ws.Range("C33").Value = 123
$ws.Range("D33").Value = 123
$ws.Range("H33").Value = 551
$ws.Range("A34").Value = 'used'
$ws.Range("B34").Value = 0.1714285714285714
$ws.Range("C34").Value = 30
$ws.Range("D34").Value = 30
$ws.Range("H34").Value = 145
$ws.Range("A35").Value = 'could'
$ws.Range("B35").Value = 0.1592356687898089
$ws.Range("C35").Value = 25
$ws.Range("D35").Value = 25
$ws.Range("H35").Value = 132
$ws.Range("A36").Value = 'work'
$ws.Range("B36").Value = 0.1582278481012658
$ws.Range("C36").Value = 50
$ws.Range("D36").Value = 50
$ws.Range("H36").Value = 266
$ws.Range("A37").Value = 'product'
$ws.Range("B37").Value = 0.1387665198237885
$ws.Range("C37").Value = 63
$ws.Range("D37").Value = 63
$ws.Range("H37").Value = 391
$ws.Range("A38").Value = 'better'
$ws.Range("B38").Value = 0.1308411214953271
$ws.Range("C38").Value = 28
$ws.Range("D38").Value = 28
$ws.Range("H38").Value = 186
$ws.Range("A39").Value = 'price'
$ws.Range("B39").Value = 0.117816091954023
$ws.Range("C39").Value = 41
$ws.Range("D39").Value = 41
$ws.Range("H39").Value = 307
$ws.Range("A40").Value = '2'
$ws.Range("B40").Value = 0.1161048689138577
$ws.Range("C40").Value = 31
$ws.Range("D40").Value = 31
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = $false
$ws.Range("H40").Value = 236
$ws.Range("B41").Value = 0.09677419354838709
$ws.Range("C41").Value = 24
$ws.Range("D41").Value = 24
$ws.Range("H41").Value = 224
$ws.Range("A42").Value = 'use'
$ws.Range("B42").Value = 0.06575342465753424
$ws.Range("C42").Value = 24
$ws.Range("D42").Value = 24
$ws.Range("H42").Value = 341
$ws.Range("A43").Value = 'like'
$ws.Range("B43").Value = 0.05766062602965404
$ws.Range("C43").Value = 35
$ws.Range("D43").Value = 36
$ws.Range("E43").Value = 0.03
$ws.Range("F43").Value = 0.97
$ws.Range("G43").Value = $true
$ws.Range("H43").Value = 572
$ws.Range("A44").Value = 'little'
$ws.Range("B44").Value = 0.05122494432071269
$ws.Range("C44").Value = 23
$ws.Range("D44").Value = 23
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = $false
$ws.Range("H44").Value = 426
$ws.Range("A45").Value = 'much'
$ws.Range("B45").Value = 0.05104408352668213
$ws.Range("C45").Value = 22
$ws.Range("E45").Value = 0.12
$ws.Range("F45").Value = 0.88
$ws.Range("G45").Value = $true
$ws.Range("H45").Value = 409
